$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: duplicate the original rows 6-8 (Provincia del Elqui / Wonderfull /
# "empedrada" reading) down into new rows 9-11, preserving their original
# values verbatim. Set the date format on the new date cells first so the
# engine reuses the existing date style instead of minting a new one.
$ws.Cells.Item(9, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(10, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(11, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(9, $col).Value = $ws.Cells.Item(6, $col).Value()
    $ws.Cells.Item(10, $col).Value = $ws.Cells.Item(7, $col).Value()
    $ws.Cells.Item(11, $col).Value = $ws.Cells.Item(8, $col).Value()
}

# Step 2: update rows 6-8 in place with the newly reported values
# (Limari province, "Sin especificar" variety, new date, new volumes/prices,
# and "granel" commercialisation unit).

# Row 6
$ws.Range("D6").Value = 44644
$ws.Range("K6").Value = "Sin especificar"
$ws.Range("M6").Value = 180
$ws.Range("Q6").Value = "$/caja 15 kilos granel"
$ws.Range("R6").Value = "Provincia de Limarí"

# Row 7
$ws.Range("D7").Value = 44644
$ws.Range("K7").Value = "Sin especificar"
$ws.Range("M7").Value = 220
$ws.Range("N7").Value = 13500
$ws.Range("O7").Value = 13500
$ws.Range("P7").Value = 13500
$ws.Range("Q7").Value = "$/caja 15 kilos granel"
$ws.Range("R7").Value = "Provincia de Limarí"
$ws.Range("S7").Value = 900

# Row 8
$ws.Range("D8").Value = 44644
$ws.Range("K8").Value = "Sin especificar"
$ws.Range("M8").Value = 290
$ws.Range("Q8").Value = "$/caja 15 kilos granel"
$ws.Range("R8").Value = "Provincia de Limarí"
